# Atualização de bases das ligas, do dia: 17-02-2024 às 11:11
#
# For a set of row pairs, the two rows' match-record data (id + all stat
# columns, i.e. column B and columns F through AC) were swapped with each
# other. Columns A (row index), C, D and E (Div / Div Original Name / Date)
# are identical between the two rows in a pair and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns whose values get swapped between the two rows of each pair.
$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

# Row pairs (1-based worksheet rows) whose content is exchanged.
$rowPairs = @(
    @(19,20),
    @(23,24),
    @(27,28),
    @(35,36),
    @(42,43),
    @(53,54),
    @(69,70),
    @(79,80),
    @(84,85),
    @(134,135),
    @(139,140),
    @(146,147),
    @(149,150),
    @(154,155),
    @(172,173),
    @(175,176),
    @(184,185),
    @(212,213),
    @(227,228),
    @(229,230)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
